$wb = $excel.ActiveWorkbook

# The workbook has two sheets that hold the same "南宁" event rows:
#   "展览"    - sheet with just the exhibition rows
#   "全部类型" - sheet aggregating all event types (same first rows)
# Both need column F ("想去人数" / interested-count) updated for rows 2-5.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 371
    $ws.Range("F3").Value = 70
    $ws.Range("F4").Value = 287
    $ws.Range("F5").Value = 4207
}
